# Generate Report for Handoff
#
# Replaces the two tracked localization files (87e26f2f...md and
# bdbeb8fd...md) with a newly generated file (d95b9b33...md) plus a
# duplicate-content handoff (ffff54c0096d...md), and flips every row's
# status from "Handed back: in sync with en-US" to "Ready for handoff"
# with refreshed handoff timestamps across the Overview / zh-cn / de-de
# sheets.

$wb = $excel.ActiveWorkbook

$oldFile1 = "87e26f2f-6832-4109-8dad-4940bb52adef.md"
$oldFile2 = "bdbeb8fd-bd65-4b25-b8e5-eaa7a2691331.md"
$newFile1 = "d95b9b33-efdc-4345-8529-9e7af8145c51.md"
$newFile2 = "ffff54c0096d-d876-4dcf-acd2-7d4a106775e3.md"

$newStatus = "Ready for handoff"
$newHoDate = "2016-08-27 12:58:56"

$newZhXlf = "d95b9b33-efdc-4345-8529-9e7af8145c51.2218bf29692761f0107d93c66046703e6ebbe90d.zh-cn.xlf"
$newDeXlf = "d95b9b33-efdc-4345-8529-9e7af8145c51.2218bf29692761f0107d93c66046703e6ebbe90d.de-de.xlf"
$newHandoffDate = "2016-08-27 12:58:50"
$newHandoffDateDe = "2016-08-27 12:58:56"
$emptyDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cc25fc090a5e4455bf5f9b6b51a5da56247956c/e2e/$newFile1", "", "", "e2e\$newFile1")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cc25fc090a5e4455bf5f9b6b51a5da56247956c/e2e/$newFile2", "", "", "e2e\$newFile2")

$wsOverview.Columns.Item(5).ColumnWidth = 16.382654825846367
$wsOverview.Columns.Item(6).ColumnWidth = 16.382654825846367

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newHandoffDate
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $emptyDate

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = $newHandoffDate
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = $emptyDate

$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("I3").Style = "Normal"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cc25fc090a5e4455bf5f9b6b51a5da56247956c/e2e/$newFile1", "", "", "$newFile1")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cc25fc090a5e4455bf5f9b6b51a5da56247956c/e2e/$newFile2", "", "", "$newFile2")

$wsZh.Columns.Item(3).ColumnWidth = 16.382654825846367
$wsZh.Columns.Item(9).ColumnWidth = 17.817272004627068
$wsZh.Columns.Item(10).ColumnWidth = 20.872143700009268

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHandoffDateDe
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $emptyDate

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = $newHandoffDateDe
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = $emptyDate

$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("I3").Style = "Normal"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b449af88b8bc8f0abec37899554346b7480ae44a/e2e/$newFile1", "", "", "$newFile1")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b449af88b8bc8f0abec37899554346b7480ae44a/e2e/$newFile2", "", "", "$newFile2")

$wsDe.Columns.Item(3).ColumnWidth = 16.382654825846367
$wsDe.Columns.Item(9).ColumnWidth = 17.817272004627068
$wsDe.Columns.Item(10).ColumnWidth = 20.872143700009268
